$d = $word.ActiveDocument

# 1. Name casing
$d.Content.Find.Execute("DHEERAJ CHAND", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dheeraj Chand", 2)

# 2. Professional title -> placeholder
$d.Content.Find.Execute("Director of Research and Analysis", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Professional Title", 2)

# 3. Phone/email formatting
$d.Content.Find.Execute("(202) 550-7110 | Dheeraj.Chand@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "202.550.7110 | dheeraj.chand@gmail.com", 2)

# 4. Years of experience
$d.Content.Find.Execute("20+ years of experience", $true, $false, $false, $false, $false,
                         $true, 1, $false, "21 years of experience", 2)

# 5. Company name redaction
$d.Content.Find.Execute("Siege Analytics, Austin, TX | 2005", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Your Company Name, Your City, ST | 2005", 2)

function Find-ParagraphIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

# 6. Remove the whole DATA PRODUCTS MANAGER / ANALYTICS SUPERVISOR /
#    SOFTWARE ENGINEER / RESEARCH DIRECTOR job blocks - everything from
#    "DATA PRODUCTS MANAGER" up to (not including) "KEY ACHIEVEMENTS AND IMPACT"
$startIdx = Find-ParagraphIndex $d "DATA PRODUCTS MANAGER"
$endIdx = (Find-ParagraphIndex $d "KEY ACHIEVEMENTS AND IMPACT") - 1
$startPar = $d.Paragraphs.Item($startIdx)
$endPar = $d.Paragraphs.Item($endIdx)
$rng = $d.Range($startPar.Range.Start, $endPar.Range.End)
$rng.Delete()

# 7. Remove the "Systems and Infrastructure Development" and
#    "Community and Stakeholder Engagement" subsections through the end
#    of the document.
$startIdx2 = Find-ParagraphIndex $d "Systems and Infrastructure Development"
$endIdx2 = $d.Paragraphs.Count
$startPar2 = $d.Paragraphs.Item($startIdx2)
$endPar2 = $d.Paragraphs.Item($endIdx2)
$rng2 = $d.Range($startPar2.Range.Start, $endPar2.Range.End)
$rng2.Delete()
